# Insert a new weekly price record for "Zapallo / Camote" at row 389.
# This pushes the existing rows 389-411 down to 390-412 (data unchanged),
# and extends the used range to A1:R412.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(389).Insert()

$ws.Range("A389").Value = 5
$ws.Range("B389").Value = "Macroferia Regional de Talca"
$ws.Range("C389").Value = "Maule"
$ws.Range("D389").Value = 45041
$ws.Range("E389").Value = 7
$ws.Range("F389").Value = 100112045
$ws.Range("G389").Value = "Zapallo"
$ws.Range("H389").Value = "Camote"
$ws.Range("I389").Value = "1a (guarda)"
$ws.Range("J389").Value = 900
$ws.Range("K389").Value = 250
$ws.Range("L389").Value = 250
$ws.Range("M389").Value = 250
$ws.Range("N389").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O389").Value = "Región del Maule"
$ws.Range("P389").Value = 250
$ws.Range("Q389").Value = 1
$ws.Range("R389").Value = "Hortaliza"
